$wb = $excel.ActiveWorkbook

# Add the new "doFindStores" sheet after the existing "doSearch" sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "doFindStores"

# Header + zip code values.
$ws.Range("A1").Value = "Zip Code"
$ws.Range("A2").Value = 19107
$ws.Range("A3").Value = 19130
$ws.Range("A4").Value = 19104
$ws.Range("A5").Value = 19152

# Match the recorded column width / selection from the target workbook.
$ws.Columns.Item(1).ColumnWidth = 13.7213541666666
$ws.Range("A2").Select() | Out-Null
